$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repair")

# --- Insert "Mean Added Vertices Percentage" above current row 7 ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Mean Added Vertices Percentage"
$ws.Range("B7").Value = 0.004448519439807632

$a7 = $ws.Range("A7")
$a7.Font.Bold = $true
$a7.HorizontalAlignment = -4108
$a7.VerticalAlignment = -4160
$a7.Borders.LineStyle = 1
$a7.Borders.Weight = 2

$b7 = $ws.Range("B7")
$b7.NumberFormat = "0.00%"

# --- Insert "Mean Added Faces Percentage" above (now) row 12 ---
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "Mean Added Faces Percentage"
$ws.Range("B12").Value = 0.001385859227250776

$a12 = $ws.Range("A12")
$a12.Font.Bold = $true
$a12.HorizontalAlignment = -4108
$a12.VerticalAlignment = -4160
$a12.Borders.LineStyle = 1
$a12.Borders.Weight = 2

$b12 = $ws.Range("B12")
$b12.NumberFormat = "0.00%"
